$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 10:16"

# Country name (column A) changes caused by re-sorting after the data refresh
$countryUpdates = @{
    24 = "Malasia"
    25 = "Dinamarca"
    26 = "Turquia"
    39 = "Indonesia"
    40 = "Islandia"
    44 = "Filipinas"
    45 = "India"
    46 = "Singapur"
    47 = "Catar"
    51 = "Hong Kong"
    52 = "Barein"
    53 = "Estonia"
    54 = "Mexico"
    55 = "Egipto"
    56 = "Croacia"
    78 = "Bosnia y Herzegovina"
    79 = "Marruecos"
    121 = "Guam"
    122 = "Montenegro"
    130 = "Guayana Francesa"
    131 = "Jamaica"
    133 = "Guatemala"
    135 = "Barbados"
    137 = "Madagascar"
}

foreach ($row in $countryUpdates.Keys) {
    $ws.Cells.Item([int]$row, 1).Value = $countryUpdates[$row]
}

# Numeric (columns B-H) updates for the refreshed figures
$valueUpdates = @{
    6 = @{ "B" = 46148; "C" = 2414; "E" = 45271 }
    15 = @{ "B" = 4628; "C" = 154; "E" = 4594 }
    22 = @{ "B" = 1960; "C" = 36; "E" = 1924 }
    24 = @{ "B" = 1624; "C" = 106; "D" = 183; "E" = 1426; "F" = 57; "G" = 1; "H" = 15 }
    25 = @{ "B" = 1577; "C" = 117; "D" = 1; "E" = 1552; "F" = 55; "G" = 0; "H" = 24 }
    26 = @{ "B" = 1529; "D" = 0; "E" = 1492; "F" = 0; "H" = 37 }
    39 = @{ "B" = 686; "C" = 107; "D" = 30; "E" = 601; "F" = 0; "G" = 6; "H" = 55 }
    40 = @{ "B" = 588; "D" = 51; "E" = 535; "F" = 14; "G" = 1; "H" = 2 }
    44 = @{ "B" = 552; "C" = 90; "D" = 19; "E" = 500; "F" = 1; "H" = 33 }
    45 = @{ "B" = 511; "C" = 12; "D" = 37; "E" = 464; "F" = 0; "H" = 10 }
    46 = @{ "B" = 509; "D" = 152; "E" = 355; "F" = 14; "H" = 2 }
    47 = @{ "C" = 0; "D" = 37; "E" = 464; "F" = 6; "H" = 0 }
    51 = @{ "B" = 386; "C" = 29; "D" = 102; "E" = 280; "F" = 4; "H" = 4 }
    52 = @{ "B" = 377; "C" = 0; "D" = 164; "E" = 211; "F" = 3; "H" = 2 }
    53 = @{ "B" = 369; "C" = 17; "E" = 365; "F" = 4; "G" = 0; "H" = 0 }
    54 = @{ "B" = 367; "C" = 51; "D" = 4; "E" = 359; "F" = 1; "G" = 1; "H" = 4 }
    55 = @{ "B" = 366; "C" = 0; "D" = 68; "E" = 279; "F" = 0; "H" = 19 }
    56 = @{ "B" = 361; "C" = 46; "D" = 5; "E" = 355; "F" = 6; "H" = 1 }
    67 = @{ "B" = 204; "C" = 18; "E" = 197 }
    78 = @{ "B" = 150; "C" = 14; "D" = 2; "E" = 147; "H" = 1 }
    79 = @{ "B" = 143; "C" = 0; "D" = 5; "H" = 4 }
    121 = @{ "C" = 0 }
    122 = @{ "C" = 2 }
    130 = @{ "B" = 23; "C" = 3; "D" = 6; "E" = 17; "H" = 0 }
    131 = @{ "B" = 21; "C" = 2; "D" = 2; "E" = 18 }
    133 = @{ "D" = 0; "E" = 19; "H" = 1 }
    135 = @{ "C" = 0 }
    137 = @{ "C" = 5 }
}

$colIndex = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5; "F" = 6; "G" = 7; "H" = 8 }

foreach ($row in $valueUpdates.Keys) {
    $rowData = $valueUpdates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item([int]$row, $colIndex[$col]).Value = $rowData[$col]
    }
}

Write-Host "Update complete"